$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F2").Value = 395246111
